$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values would be auto-parsed as numbers by Excel need
# an explicit Text number format so they are stored as text, matching
# the source data (which is text throughout column D/E).
$textCells = @('D5', 'D6', 'D7', 'D10', 'D11', 'D12', 'D13', 'D14', 'D20', 'D21', 'D22', 'D24', 'D30', 'D34', 'D36', 'D40', 'D41', 'D45', 'D48', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values.
$ws.Range('D2').Value = '65.950.32'
$ws.Range('E2').Value = '  +6.94%  '
$ws.Range('D3').Value = '3.016.20'
$ws.Range('E3').Value = '  +4.13%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '585.75'
$ws.Range('E5').Value = '  +3.20%  '
$ws.Range('D6').Value = '156.43'
$ws.Range('E6').Value = '  +9.23%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.012.76'
$ws.Range('E8').Value = '  +4.07%  '
$ws.Range('E9').Value = '  +3.20%  '
$ws.Range('D10').Value = '7.06'
$ws.Range('E10').Value = '  +2.07%  '
$ws.Range('D11').Value = '0.156'
$ws.Range('E11').Value = '  +7.05%  '
$ws.Range('D12').Value = '0.452'
$ws.Range('E12').Value = '  +5.28%  '
$ws.Range('D13').Value = '0.0000252'
$ws.Range('E13').Value = '  +9.17%  '
$ws.Range('D14').Value = '34.52'
$ws.Range('E14').Value = '  +8.55%  '
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').Value = '65.913.35'
$ws.Range('E16').Value = '  +7.01%  '
$ws.Range('D17').Value = '3.515.55'
$ws.Range('E17').Value = '  +4.12%  '
$ws.Range('E18').Value = '  +6.31%  '
$ws.Range('D19').Value = '3.016.00'
$ws.Range('E19').Value = '  +4.14%  '
$ws.Range('D20').Value = '465.00'
$ws.Range('E20').Value = '  +7.67%  '
$ws.Range('D21').Value = '13.86'
$ws.Range('E21').Value = '  +6.30%  '
$ws.Range('D22').Value = '0.684'
$ws.Range('E22').Value = '  +4.67%  '
$ws.Range('E23').Value = '  +8.28%  '
$ws.Range('D24').Value = '82.25'
$ws.Range('E24').Value = '  +3.90%  '
$ws.Range('E25').Value = '  +5.27%  '
$ws.Range('E26').Value = '  +11.82%  '
$ws.Range('E27').Value = '  +8.02%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +13.46%  '
$ws.Range('D30').Value = '2.40'
$ws.Range('E30').Value = '  +17.67%  '
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('E32').Value = '  +5.05%  '
$ws.Range('E33').Value = '  +5.28%  '
$ws.Range('D34').Value = '27.02'
$ws.Range('E34').Value = '  +5.92%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +4.39%  '
$ws.Range('E37').Value = '  +8.38%  '
$ws.Range('E38').Value = '  +13.07%  '
$ws.Range('E39').Value = '  +7.57%  '
$ws.Range('D40').Value = '49.30'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('D41').Value = '44.69'
$ws.Range('E41').Value = '  +12.28%  '
$ws.Range('E42').Value = '  +8.61%  '
$ws.Range('E43').Value = '  +12.21%  '
$ws.Range('E44').Value = '  +3.73%  '
$ws.Range('D45').Value = '395.31'
$ws.Range('E45').Value = '  +15.06%  '
$ws.Range('D46').Value = '2.805.65'
$ws.Range('E46').Value = '  +4.35%  '
$ws.Range('E47').Value = '  +5.82%  '
$ws.Range('D48').Value = '134.16'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').Value = '23.83'
$ws.Range('E50').Value = '  +10.66%  '
$ws.Range('E51').Value = '  +4.23%  '
